$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$headers = $sec.Headers
$hdr = $headers.Item(1)
$r = $hdr.Range
$tbl = $r.Tables.Item(1)
$cell1 = $tbl.Cell(1,1)
$cell1.Range.Font.Italic = 1
Write-Host "done"
